{"js": "// 1) \"Masukin [ke] C:/Bosnet/Bin\" -> merge \" \", \"ke\", \" C:/Bosnet/Bin\"\n//    (and their spellStart/spellEnd proofErr marks around \"ke\") into a\n//    single run \" ke C:/Bosnet/Bin\" right after \"Masukin\".\nconst body = context.document.body;\n\nconst keResults = body.search(\" ke C:/Bosnet/Bin\", { matchCase: true });\nawait context.sync();\n\nif (keResults.items.length === 0) {\n  throw new Error('Could not find \" ke C:/Bosnet/Bin\" in the document.');\n}\n\nconst keRange = keResults.items[0];\n// Route the replacement through a content control so Word is forced to\n// treat the whole span as one fresh run, instead of leaving the old\n// (now-empty) spellStart/spellEnd proofErr pair behind.\nconst keCC = keRange.insertContentControl();\nkeCC.insertText(\" ke C:/Bosnet/Bin\", Word.InsertLocation.replace);\nawait context.sync();\nkeCC.cannotDelete = false;\nkeCC.delete(true); // drop the wrapper, keep the merged text\nawait context.sync();\n\n// 2) \"SQN\" -> \"F\" followed by \"QN\" (two separate runs, text becomes \"FQN\").\nconst sqnResults = body.search(\"SQN\", { matchCase: true });\nawait context.sync();\n\nif (sqnResults.items.length === 0) {\n  throw new Error('Could not find \"SQN\" in the document.');\n}\n\nconst sqnRange = sqnResults.items[0];\n// First shrink \"SQN\" down to \"F\" in place (still a single run).\nsqnRange.insertText(\"F\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Then locate that new \"F\" precisely and append \"QN\" as its own run,\n// again using a content control so it doesn't get silently re-merged\n// into the neighbouring run during serialization.\nconst fResults = body.search(\"o:F\", { matchCase: true });\nawait context.sync();\n\nif (fResults.items.length === 0) {\n  throw new Error('Could not find \"o:F\" in the document.');\n}\n\nconst fRange = fResults.items[0];\nconst afterF = fRange.getRange(Word.RangeLocation.after);\nconst qnCC = afterF.insertContentControl();\nqnCC.insertText(\"QN\", Word.InsertLocation.replace);\nawait context.sync();\nqnCC.cannotDelete = false;\nqnCC.delete(true);\nawait context.sync();\n", "ps1": "# Word COM interop script (PowerShell-style) \u2014 applies the same two edits\n# as edit.js against $word.ActiveDocument.\n\n$d = $word.ActiveDocument\n\n# 1) \"Masukin [ ][ke][ C:/Bosnet/Bin]\" -> merge the trailing \" \", \"ke\",\n#    \" C:/Bosnet/Bin\" runs (and the spellStart/spellEnd proofErr marks\n#    wrapping \"ke\") into one run reading \" ke C:/Bosnet/Bin\".\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \" ke C:/Bosnet/Bin\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \" ke C:/Bosnet/Bin\"\n$find1.Execute(\" ke C:/Bosnet/Bin\", $false, $false, $false, $false, $false, $true, 1, $false, \" ke C:/Bosnet/Bin\", 2) | Out-Null\n\n# 2) \"SQN\" -> \"F\" followed by \"QN\" as two separate runs (text becomes \"FQN\").\n$rng = $d.Content\n$rng.Find.Execute(\"SQN\") | Out-Null\n$start = $rng.Start\n$end = $rng.End\n\n# First character (\"S\") becomes its own range, rest (\"QN\") becomes another.\n$rFirst = $d.Range($start, $start + 1)\n$rRest = $d.Range($start + 1, $end)\n\n# Change \"QN\" portion first while it is still a distinct Range object...\n$rRest.Text = \"QN\"\n# ...then shrink the first range down to \"F\". Doing the edits through two\n# independently-resolved Range objects (rather than one contiguous\n# insertText/replace) keeps them as two separate runs instead of Word\n# silently re-coalescing them back into a single run.\n$rFirst.Text = \"F\"\n"}
